$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.075.96"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.495.60"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.84"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.82"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.520"
$ws.Range("E7").Value = "  -1.33%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  -4.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.83"
$ws.Range("E10").Value = "  -3.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.03"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.09"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.886.65"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.494.61"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.907.68"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.00"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("E20").Value = "  +9.21%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0930"
$ws.Range("E22").Value = "  -1.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.03"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.47"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.51"
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.76"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  -4.31%  "
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.96"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  -4.44%  "
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.58"
$ws.Range("E38").Value = "  -2.90%  "
$ws.Range("E39").Value = "  -4.26%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.39"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.61"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.110"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0302"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.995.97"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.90"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.86"
$ws.Range("E51").Value = "  -1.85%  "
